# Apply the AP Invoices Import row-2 data corrections described in the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# GL_Date: 08/04/25 -> 08/06/25 (keep as text, not an auto-converted date serial)
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "08/06/25"

# Vendor_Code: LORSON -> JONSUP
$ws.Range("E2").Value = "JONSUP"

# Invoice_Number: 01-871062 -> 110-S10112669.001
$ws.Range("G2").Value = "110-S10112669.001"

# Invoice_Date: 08/04/25 -> 08/06/25 (keep as text, not an auto-converted date serial)
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "08/06/25"

# Invoice_Amount: 89 -> 99.31
$ws.Range("J2").Value = 99.31

# Remarks: (blank) -> San Leandro Unit Install
$ws.Range("R2").Value = "San Leandro Unit Install"

# Distribution_GL_Account: 5030 -> 5260
$ws.Range("T2").Value = 5260

# Job_Number: 25.09 -> (blank)
$ws.Range("U2").Value = ""

# Phase_Code: 320 -> (removed)
$ws.Range("V2").ClearContents()

# Cost_Type: M -> (removed)
$ws.Range("W2").ClearContents()

# WO_Number: (blank) -> 2025 (kept as text)
$ws.Range("AA2").NumberFormat = "@"
$ws.Range("AA2").Value = "2025"

# Item_Code: !Material -> !Service Material
$ws.Range("AB2").Value = "!Service Material"
